$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.76548733449558
$ws.Range("D2").Value = 3.332866452911979
$ws.Range("E2").Value = 16.7405582290129
$ws.Range("F2").Value = 18.35869918359322
$ws.Range("G2").Value = 21.97987378217738
$ws.Range("H2").Value = 10.62994996385521
$ws.Range("N2").Value = 18.20530214750648
$ws.Range("O2").Value = 15.51357751095988

$ws.Range("B3").Value = 7.690488708080919
$ws.Range("D3").Value = 3.331906116542084
$ws.Range("E3").Value = 15.77358072534926
$ws.Range("F3").Value = 17.85750361967092
$ws.Range("G3").Value = 20.88400046637141
$ws.Range("H3").Value = 10.56398091616973
$ws.Range("N3").Value = 17.93039231733595
$ws.Range("O3").Value = 15.21194444154728

$ws.Range("B4").Value = 7.645832371061269
$ws.Range("D4").Value = 3.331667210034011
$ws.Range("E4").Value = 15.15375536400619
$ws.Range("F4").Value = 17.54855801098144
$ws.Range("G4").Value = 20.18742174270481
$ws.Range("H4").Value = 10.52596340044021
$ws.Range("N4").Value = 17.76180298315366
$ws.Range("O4").Value = 15.02856969616999

$ws.Range("B5").Value = 7.62800595514101
$ws.Range("D5").Value = 3.33165832644058
$ws.Range("E5").Value = 14.89487218464948
$ws.Range("F5").Value = 17.4225731203034
$ws.Range("G5").Value = 19.89801699846327
$ws.Range("H5").Value = 10.51110920613297
$ws.Range("N5").Value = 17.69323820309472
$ws.Range("O5").Value = 14.9544145422213

$ws.Range("B6").Value = 7.625068927627218
$ws.Range("D6").Value = 3.331662202569909
$ws.Range("E6").Value = 14.85151310981644
$ws.Range("F6").Value = 17.40165469098864
$ws.Range("G6").Value = 19.84964027438417
$ws.Range("H6").Value = 10.50868157673871
$ws.Range("N6").Value = 17.68186382811895
$ws.Range("O6").Value = 14.94213892303543

$ws.Range("B7").Value = 7.645590426598107
$ws.Range("D7").Value = 3.331666731627296
$ws.Range("E7").Value = 15.1502890936753
$ws.Range("F7").Value = 17.54685899359088
$ws.Range("G7").Value = 20.18354057289569
$ws.Range("H7").Value = 10.52576047143519
$ws.Range("N7").Value = 17.76087762817026
$ws.Range("O7").Value = 15.02756715204569

$ws.Range("B8").Value = 7.739351513268095
$ws.Range("D8").Value = 3.332462701585339
$ws.Range("E8").Value = 16.41270293448153
$ws.Range("F8").Value = 18.18626524397519
$ws.Range("G8").Value = 21.60716863689971
$ws.Range("H8").Value = 10.60669544107147
$ws.Range("N8").Value = 18.11052099476146
$ws.Range("O8").Value = 15.4092615210346

$ws.Range("B9").Value = 7.9333151054759
$ws.Range("D9").Value = 3.336794265068324
$ws.Range("E9").Value = 18.83526645892429
$ws.Range("F9").Value = 19.42152583551522
$ws.Range("G9").Value = 24.19549305282739
$ws.Range("H9").Value = 10.78454097625027
$ws.Range("N9").Value = 18.79413998554022
$ws.Range("O9").Value = 16.16742811790309

$ws.Range("B10").Value = 8.080666255167973
$ws.Range("D10").Value = 3.341648870266802
$ws.Range("E10").Value = 20.51580639411464
$ws.Range("F10").Value = 20.30649208104569
$ws.Range("G10").Value = 25.95640391697248
$ws.Range("H10").Value = 10.92598972804872
$ws.Range("N10").Value = 19.2903645072747
$ws.Range("O10").Value = 16.72400890826519

$ws.Range("B11").Value = 8.14846916618764
$ws.Range("D11").Value = 3.344216520351163
$ws.Range("E11").Value = 21.23749084471423
$ws.Range("F11").Value = 20.70207008789831
$ws.Range("G11").Value = 26.72430725654529
$ws.Range("H11").Value = 10.99247610726923
$ws.Range("N11").Value = 19.51383106981896
$ws.Range("O11").Value = 16.97584131324049

$ws.Range("B12").Value = 8.17423166240809
$ws.Range("D12").Value = 3.345240104828057
$ws.Range("E12").Value = 21.50465745367265
$ws.Range("F12").Value = 20.85070168403498
$ws.Range("G12").Value = 27.01015662864143
$ws.Range("H12").Value = 11.01794236505454
$ws.Range("N12").Value = 19.59804885977592
$ws.Range("O12").Value = 17.07090818264913

$ws.Range("B13").Value = 8.168679758780518
$ws.Range("D13").Value = 3.345017383736601
$ws.Range("E13").Value = 21.44738995832763
$ws.Range("F13").Value = 20.81874541673822
$ws.Range("G13").Value = 26.94881594057337
$ws.Range("H13").Value = 11.01244518410014
$ws.Range("N13").Value = 19.57993023154881
$ws.Range("O13").Value = 17.05044857741606

$ws.Range("B14").Value = 8.150587055476443
$ws.Range("D14").Value = 3.344299706018385
$ws.Range("E14").Value = 21.25959324810564
$ws.Range("F14").Value = 20.71432231239457
$ws.Range("G14").Value = 26.74792412138484
$ws.Range("H14").Value = 10.99456554138558
$ws.Range("N14").Value = 19.52076825227033
$ws.Range("O14").Value = 16.98366902686251

$ws.Range("B15").Value = 8.139515373066374
$ws.Range("D15").Value = 3.343866772780879
$ws.Range("E15").Value = 21.14376642488964
$ws.Range("F15").Value = 20.65020392082526
$ws.Range("G15").Value = 26.62422428068378
$ws.Range("H15").Value = 10.98365086186756
$ws.Range("N15").Value = 19.48447494387068
$ws.Range("O15").Value = 16.94272308504471

$ws.Range("B16").Value = 8.076248655077
$ws.Range("D16").Value = 3.341488264673342
$ws.Range("E16").Value = 20.46778435323411
$ws.Range("F16").Value = 20.28048541154647
$ws.Range("G16").Value = 25.9055371870232
$ws.Range("H16").Value = 10.921686192136
$ws.Range("N16").Value = 19.27570805562915
$ws.Range("O16").Value = 16.70751495472316

$ws.Range("B17").Value = 8.037617191288184
$ws.Range("D17").Value = 3.340120871658005
$ws.Range("E17").Value = 20.04215383634661
$ws.Range("F17").Value = 20.05176529754235
$ws.Range("G17").Value = 25.45602378938278
$ws.Range("H17").Value = 10.88420815208812
$ws.Range("N17").Value = 19.14699881872644
$ws.Range("O17").Value = 16.56280069720605

$ws.Range("B18").Value = 8.01547150347373
$ws.Range("D18").Value = 3.339368223788723
$ws.Range("E18").Value = 19.79331164051828
$ws.Range("F18").Value = 19.91956150706805
$ws.Range("G18").Value = 25.19436406922753
$ws.Range("H18").Value = 10.86285440835427
$ws.Range("N18").Value = 19.07275961144546
$ws.Range("O18").Value = 16.47944270756341

$ws.Range("B19").Value = 8.00798681698601
$ws.Range("D19").Value = 3.339119213877821
$ws.Range("E19").Value = 19.70836458320055
$ws.Range("F19").Value = 19.87469300417116
$ws.Range("G19").Value = 25.10524197005105
$ws.Range("H19").Value = 10.85565975489246
$ws.Range("N19").Value = 19.04759001474569
$ws.Range("O19").Value = 16.45120130807218

$ws.Range("B20").Value = 8.041722081226251
$ws.Range("D20").Value = 3.340262933356617
$ws.Range("E20").Value = 20.08787984548762
$ws.Range("F20").Value = 20.07618141573258
$ws.Range("G20").Value = 25.50419863862297
$ws.Range("H20").Value = 10.88817691407267
$ws.Range("N20").Value = 19.16072234871258
$ws.Range("O20").Value = 16.57821917134406

$ws.Range("B21").Value = 8.155899150234847
$ws.Range("D21").Value = 3.344509117089792
$ws.Range("E21").Value = 21.31491949976721
$ws.Range("F21").Value = 20.74502673075174
$ws.Range("G21").Value = 26.80706618550222
$ws.Range("H21").Value = 10.99980952428328
$ws.Range("N21").Value = 19.53815713259456
$ws.Range("O21").Value = 17.00329262939363

$ws.Range("B22").Value = 8.23101661654694
$ws.Range("D22").Value = 3.347582896093313
$ws.Range("E22").Value = 22.08120714169245
$ws.Range("F22").Value = 21.17530009767951
$ws.Range("G22").Value = 27.62971804028947
$ws.Range("H22").Value = 11.07444526050048
$ws.Range("N22").Value = 19.78244317294325
$ws.Range("O22").Value = 17.27933077334015

$ws.Range("B23").Value = 8.190887495043029
$ws.Range("D23").Value = 3.345915169710971
$ws.Range("E23").Value = 21.67547574966455
$ws.Range("F23").Value = 20.94633133725022
$ws.Range("G23").Value = 27.1933421629898
$ws.Range("H23").Value = 11.03446353276845
$ws.Range("N23").Value = 19.65230661567175
$ws.Range("O23").Value = 17.13219820154773

$ws.Range("B24").Value = 8.039866059139561
$ws.Range("D24").Value = 3.340198602956081
$ws.Range("E24").Value = 20.06722001980257
$ws.Range("F24").Value = 20.06514508651015
$ws.Range("G24").Value = 25.48242883184161
$ws.Range("H24").Value = 10.8863820345403
$ws.Range("N24").Value = 19.15451869034855
$ws.Range("O24").Value = 16.57124896922864

$ws.Range("B25").Value = 7.879895846024697
$ws.Range("D25").Value = 3.33532750643249
$ws.Range("E25").Value = 18.1782030260668
$ws.Range("F25").Value = 19.09055749939816
$ws.Range("G25").Value = 23.51897299350291
$ws.Range("H25").Value = 10.73446794358163
$ws.Range("N25").Value = 18.60992361498234
$ws.Range("O25").Value = 15.96195547730649
